# Init Realtime prices tab
# Refresh the crawler results (B:G) for all 22 tickers (rows 2-23) with the
# latest snapshot taken on 2023-05-17 17:50:07.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 10
$ws.Cells.Item(2, 3).Value = 3
$ws.Cells.Item(2, 4).Value = 23
$ws.Cells.Item(2, 5).Value = 9.300000000000001
$ws.Cells.Item(2, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(2, 7).Value = "rising"

$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 13
$ws.Cells.Item(3, 5).Value = 9.300000000000001
$ws.Cells.Item(3, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(3, 7).Value = "rising"

$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(4, 7).Value = "falling"

$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 7
$ws.Cells.Item(5, 5).Value = 3.7
$ws.Cells.Item(5, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(5, 7).Value = "rising"

$ws.Cells.Item(6, 2).Value = 4
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 8
$ws.Cells.Item(6, 5).Value = 3.8
$ws.Cells.Item(6, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(6, 7).Value = "rising"

$ws.Cells.Item(7, 2).Value = 8
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 8
$ws.Cells.Item(7, 5).Value = 7.8
$ws.Cells.Item(7, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(7, 7).Value = "rising"

$ws.Cells.Item(8, 2).Value = 5
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 12
$ws.Cells.Item(8, 5).Value = 5.2
$ws.Cells.Item(8, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(8, 7).Value = "stable"

$ws.Cells.Item(9, 2).Value = 3
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 4).Value = 3
$ws.Cells.Item(9, 5).Value = 1.3
$ws.Cells.Item(9, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(9, 7).Value = "rising"

$ws.Cells.Item(10, 2).Value = 2
$ws.Cells.Item(10, 3).Value = 1
$ws.Cells.Item(10, 4).Value = 2
$ws.Cells.Item(10, 5).Value = 1.2
$ws.Cells.Item(10, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(10, 7).Value = "rising"

$ws.Cells.Item(11, 2).Value = 10
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 4).Value = 13
$ws.Cells.Item(11, 5).Value = 9.300000000000001
$ws.Cells.Item(11, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(11, 7).Value = "rising"

$ws.Cells.Item(12, 2).Value = 14
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = 15
$ws.Cells.Item(12, 5).Value = 13.5
$ws.Cells.Item(12, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(12, 7).Value = "rising"

$ws.Cells.Item(13, 2).Value = 12
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 4).Value = 18
$ws.Cells.Item(13, 5).Value = 11.8
$ws.Cells.Item(13, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(13, 7).Value = "rising"

$ws.Cells.Item(14, 2).Value = 15
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 4).Value = 26
$ws.Cells.Item(14, 5).Value = 15.6
$ws.Cells.Item(14, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(14, 7).Value = "rising"

$ws.Cells.Item(15, 2).Value = 6
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(15, 4).Value = 7
$ws.Cells.Item(15, 5).Value = 3.7
$ws.Cells.Item(15, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(15, 7).Value = "rising"

$ws.Cells.Item(16, 2).Value = 6
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 23
$ws.Cells.Item(16, 5).Value = 8.300000000000001
$ws.Cells.Item(16, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(16, 7).Value = "rising"

$ws.Cells.Item(17, 2).Value = 11
$ws.Cells.Item(17, 3).Value = 9
$ws.Cells.Item(17, 4).Value = 24
$ws.Cells.Item(17, 5).Value = 4.4
$ws.Cells.Item(17, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(17, 7).Value = "stable"

$ws.Cells.Item(18, 2).Value = 3
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 8
$ws.Cells.Item(18, 5).Value = 3.8
$ws.Cells.Item(18, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(18, 7).Value = "rising"

$ws.Cells.Item(19, 2).Value = 14
$ws.Cells.Item(19, 3).Value = 2
$ws.Cells.Item(19, 4).Value = 24
$ws.Cells.Item(19, 5).Value = 14.4
$ws.Cells.Item(19, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(19, 7).Value = "rising"

$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 9
$ws.Cells.Item(20, 5).Value = -0.1000000000000001
$ws.Cells.Item(20, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(20, 7).Value = "falling"

$ws.Cells.Item(21, 2).Value = 5
$ws.Cells.Item(21, 3).Value = 3
$ws.Cells.Item(21, 4).Value = 13
$ws.Cells.Item(21, 5).Value = 3.3
$ws.Cells.Item(21, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(21, 7).Value = "rising"

$ws.Cells.Item(22, 2).Value = 2
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 4).Value = 19
$ws.Cells.Item(22, 5).Value = 1.9
$ws.Cells.Item(22, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(22, 7).Value = "rising"

$ws.Cells.Item(23, 2).Value = 3
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = 9
$ws.Cells.Item(23, 5).Value = 1.9
$ws.Cells.Item(23, 6).Value = "2023-05-17 17:50:07"
$ws.Cells.Item(23, 7).Value = "rising"

Write-Host "Updated rows 2-23 (B:G) with latest realtime price snapshot."
